$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new task entry
$ws.Range("D3").Value = "Imported Promotion`nOrder in UAT for 7 regions"
$ws.Range("E3").Value = "Imported Promotion`nOrder in UAT for 7 regions"
$ws.Range("F3").Value = "Centra"
$ws.Range("G3").Value = "UAT"
$ws.Range("H3").Value = "Completed"
$ws.Range("V3").Value = "Manirathnam"

# Wrap text + row height for row 3
$ws.Range("D3:E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 24.25

# Row 4: new task entry
$ws.Range("D4").Value = "Tested centra. Bugs found(2)"
$ws.Range("E4").Value = "Tested centra. Bugs found(2)"
$ws.Range("F4").Value = "Centra"
$ws.Range("G4").Value = 8098
$ws.Range("H4").Value = "Inprogress"

# E4 picked up the row-1 style (distinct xf slot) during editing
$ws.Range("E4").Style = $ws.Rows.Item(1).Style

$ws.Range("H4").Select()
